$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header while it is still at its original location (E1),
# so the shared-string entry is updated in place rather than creating a new one.
$ws.Range("E1").Value = "MODELCONDITION"

# Remove the obsolete first column (A) entirely; this shifts B:F left to A:E,
# dropping column A's old contents (2 / 9 with bold border style) and moving
# the header/data columns into their new positions.
$ws.Columns("A:A").Delete()
